$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E7").Value = 12.3409
$ws.Range("A8").Value = -21.15580000000001
$ws.Range("A10").Value = -20.50089999999997
$ws.Range("A12").Value = -22.40550000000003
$ws.Range("C13").Value = -13.4795
$ws.Range("A18").Value = -22.30330000000002
$ws.Range("E20").Value = 12.31419999999999
$ws.Range("A25").Value = -22.25990000000003
